$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K column (최종점수 / final score) values
$ws.Range("K2").Value = 54.7
$ws.Range("K3").Value = 50.7
$ws.Range("K4").Value = 46.5
$ws.Range("K5").Value = 45.3
$ws.Range("K6").Value = 36.7

# Update N column (MACRO_SCORE) values
$ws.Range("N2").Value = 51.15965480231979
$ws.Range("N3").Value = 51.15965480231979
$ws.Range("N4").Value = 51.15965480231979
$ws.Range("N5").Value = 51.15965480231979
$ws.Range("N6").Value = 51.15965480231979
